$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 133
$ws.Range("F3").Value = 1314
$ws.Range("F5").Value = 1001
$ws.Range("F6").Value = 1776
$ws.Range("F7").Value = 550
$ws.Range("F8").Value = 1178
$ws.Range("F9").Value = 54
$ws.Range("F10").Value = 11
$ws.Range("F12").Value = 280
$ws.Range("F13").Value = 62
$ws.Range("F15").Value = 665
$ws.Range("F17").Value = 99
$ws.Range("F21").Value = 135
$ws.Range("F22").Value = 662
$ws.Range("F24").Value = 638
$ws.Range("F29").Value = 153
$ws.Range("F30").Value = 38
$ws.Range("F31").Value = 264

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 119

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 306

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 306
$ws.Range("F3").Value = 133
$ws.Range("F4").Value = 1314
$ws.Range("F6").Value = 1001
$ws.Range("F7").Value = 1776
$ws.Range("F8").Value = 550
$ws.Range("F9").Value = 1178
$ws.Range("F10").Value = 54
$ws.Range("F12").Value = 11
$ws.Range("F14").Value = 280
$ws.Range("F15").Value = 62
$ws.Range("F17").Value = 665
$ws.Range("F19").Value = 99
$ws.Range("F29").Value = 135
$ws.Range("F30").Value = 662
$ws.Range("F32").Value = 638
$ws.Range("F39").Value = 153
$ws.Range("F40").Value = 38
$ws.Range("F41").Value = 264
$ws.Range("F43").Value = 119
